$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Cells.Item(2, 7).Value = 6.775549333333333
$ws.Cells.Item(2, 8).Value = 20.326648
$ws.Cells.Item(2, 9).Value = 0.5307754563424079
$ws.Cells.Item(2, 10).Value = 0.5307754563424079
$ws.Cells.Item(2, 11).Value = 3
$ws.Cells.Item(2, 12).Value = 1
$ws.Cells.Item(2, 13).Value = 14.91571833333333
$ws.Cells.Item(2, 14).Value = 44.747155
$ws.Cells.Item(2, 15).Value = 0.9626318069253016
$ws.Cells.Item(2, 16).Value = 0.9626318069253015
$ws.Cells.Item(2, 17).Value = 101.0621854096044
$ws.Cells.Item(2, 18).Value = 909.55966868644
$ws.Cells.Item(2, 19).Value = 0.5109413366104936
$ws.Cells.Item(2, 20).Value = 0.5109413366104936
# Row 3
$ws.Cells.Item(3, 7).Value = 6.775549333333333
$ws.Cells.Item(3, 8).Value = 20.326648
$ws.Cells.Item(3, 9).Value = 0.5307754563424079
$ws.Cells.Item(3, 10).Value = 0.5307754563424079
$ws.Cells.Item(3, 15).Value = 0.01362824797293961
$ws.Cells.Item(3, 16).Value = 0.01362824797293961
$ws.Cells.Item(3, 17).Value = 1.430765650522667
$ws.Cells.Item(3, 18).Value = 12.876890854704
$ws.Cells.Item(3, 19).Value = 0.007233539536984518
$ws.Cells.Item(3, 20).Value = 0.007233539536984517
# Row 4
$ws.Cells.Item(4, 7).Value = 6.775549333333333
$ws.Cells.Item(4, 8).Value = 20.326648
$ws.Cells.Item(4, 9).Value = 0.5307754563424079
$ws.Cells.Item(4, 10).Value = 0.5307754563424079
$ws.Cells.Item(4, 13).Value = 0.3355976666666667
$ws.Cells.Item(4, 14).Value = 1.006793
$ws.Cells.Item(4, 15).Value = 0.02165882869625444
$ws.Cells.Item(4, 16).Value = 0.02165882869625443
$ws.Cells.Item(4, 17).Value = 2.273858546651556
$ws.Cells.Item(4, 18).Value = 20.464726919864
$ws.Cells.Item(4, 19).Value = 0.01149597468509649
$ws.Cells.Item(4, 20).Value = 0.01149597468509649
# Row 5
$ws.Cells.Item(5, 7).Value = 6.775549333333333
$ws.Cells.Item(5, 8).Value = 20.326648
$ws.Cells.Item(5, 9).Value = 0.5307754563424079
$ws.Cells.Item(5, 10).Value = 0.5307754563424079
$ws.Cells.Item(5, 13).Value = 0.03224633333333333
$ws.Cells.Item(5, 14).Value = 0.09673899999999999
$ws.Cells.Item(5, 15).Value = 0.002081116405504366
$ws.Cells.Item(5, 16).Value = 0.002081116405504366
$ws.Cells.Item(5, 17).Value = 0.2184866223191111
$ws.Cells.Item(5, 18).Value = 1.966379600872
$ws.Cells.Item(5, 19).Value = 0.001104605509833252
$ws.Cells.Item(5, 20).Value = 0.001104605509833252
# Row 6
$ws.Cells.Item(6, 9).Value = 0.3421215311185197
$ws.Cells.Item(6, 10).Value = 0.3421215311185197
$ws.Cells.Item(6, 11).Value = 3
$ws.Cells.Item(6, 12).Value = 1
$ws.Cells.Item(6, 13).Value = 14.91571833333333
$ws.Cells.Item(6, 14).Value = 44.747155
$ws.Cells.Item(6, 15).Value = 0.9626318069253016
$ws.Cells.Item(6, 16).Value = 0.9626318069253015
$ws.Cells.Item(6, 17).Value = 65.14157577816222
$ws.Cells.Item(6, 18).Value = 586.2741820034599
$ws.Cells.Item(6, 19).Value = 0.3293370676886714
$ws.Cells.Item(6, 20).Value = 0.3293370676886714
# Row 7
$ws.Cells.Item(7, 9).Value = 0.3421215311185197
$ws.Cells.Item(7, 10).Value = 0.3421215311185197
$ws.Cells.Item(7, 15).Value = 0.01362824797293961
$ws.Cells.Item(7, 16).Value = 0.01362824797293961
$ws.Cells.Item(7, 19).Value = 0.004662517062964964
$ws.Cells.Item(7, 20).Value = 0.004662517062964963
# Row 8
$ws.Cells.Item(8, 9).Value = 0.3421215311185197
$ws.Cells.Item(8, 10).Value = 0.3421215311185197
$ws.Cells.Item(8, 13).Value = 0.3355976666666667
$ws.Cells.Item(8, 14).Value = 1.006793
$ws.Cells.Item(8, 15).Value = 0.02165882869625444
$ws.Cells.Item(8, 16).Value = 0.02165882869625443
$ws.Cells.Item(8, 17).Value = 1.465659269341778
$ws.Cells.Item(8, 18).Value = 13.190933424076
$ws.Cells.Item(8, 19).Value = 0.007409951635796301
$ws.Cells.Item(8, 20).Value = 0.0074099516357963
# Row 9
$ws.Cells.Item(9, 9).Value = 0.3421215311185197
$ws.Cells.Item(9, 10).Value = 0.3421215311185197
$ws.Cells.Item(9, 13).Value = 0.03224633333333333
$ws.Cells.Item(9, 14).Value = 0.09673899999999999
$ws.Cells.Item(9, 15).Value = 0.002081116405504366
$ws.Cells.Item(9, 16).Value = 0.002081116405504366
$ws.Cells.Item(9, 17).Value = 0.1408297555275555
$ws.Cells.Item(9, 18).Value = 1.267467799748
$ws.Cells.Item(9, 19).Value = 0.000711994731087024
$ws.Cells.Item(9, 20).Value = 0.000711994731087024
# Row 10
$ws.Cells.Item(10, 7).Value = 1.622518
$ws.Cells.Item(10, 8).Value = 4.867554
$ws.Cells.Item(10, 9).Value = 0.1271030125390725
$ws.Cells.Item(10, 10).Value = 0.1271030125390725
$ws.Cells.Item(10, 11).Value = 3
$ws.Cells.Item(10, 12).Value = 1
$ws.Cells.Item(10, 13).Value = 14.91571833333333
$ws.Cells.Item(10, 14).Value = 44.747155
$ws.Cells.Item(10, 15).Value = 0.9626318069253016
$ws.Cells.Item(10, 16).Value = 0.9626318069253015
$ws.Cells.Item(10, 17).Value = 24.20102147876333
$ws.Cells.Item(10, 18).Value = 217.80919330887
$ws.Cells.Item(10, 19).Value = 0.1223534026261366
$ws.Cells.Item(10, 20).Value = 0.1223534026261366
# Row 11
$ws.Cells.Item(11, 7).Value = 1.622518
$ws.Cells.Item(11, 8).Value = 4.867554
$ws.Cells.Item(11, 9).Value = 0.1271030125390725
$ws.Cells.Item(11, 10).Value = 0.1271030125390725
$ws.Cells.Item(11, 15).Value = 0.01362824797293961
$ws.Cells.Item(11, 16).Value = 0.01362824797293961
$ws.Cells.Item(11, 17).Value = 0.342620635988
$ws.Cells.Item(11, 18).Value = 3.083585723892
$ws.Cells.Item(11, 19).Value = 0.001732191372990133
$ws.Cells.Item(11, 20).Value = 0.001732191372990133
# Row 12
$ws.Cells.Item(12, 7).Value = 1.622518
$ws.Cells.Item(12, 8).Value = 4.867554
$ws.Cells.Item(12, 9).Value = 0.1271030125390725
$ws.Cells.Item(12, 10).Value = 0.1271030125390725
$ws.Cells.Item(12, 13).Value = 0.3355976666666667
$ws.Cells.Item(12, 14).Value = 1.006793
$ws.Cells.Item(12, 15).Value = 0.02165882869625444
$ws.Cells.Item(12, 16).Value = 0.02165882869625443
$ws.Cells.Item(12, 17).Value = 0.5445132549246667
$ws.Cells.Item(12, 18).Value = 4.900619294322
$ws.Cells.Item(12, 19).Value = 0.002752902375361651
$ws.Cells.Item(12, 20).Value = 0.00275290237536165
# Row 13
$ws.Cells.Item(13, 7).Value = 1.622518
$ws.Cells.Item(13, 8).Value = 4.867554
$ws.Cells.Item(13, 9).Value = 0.1271030125390725
$ws.Cells.Item(13, 10).Value = 0.1271030125390725
$ws.Cells.Item(13, 13).Value = 0.03224633333333333
$ws.Cells.Item(13, 14).Value = 0.09673899999999999
$ws.Cells.Item(13, 15).Value = 0.002081116405504366
$ws.Cells.Item(13, 16).Value = 0.002081116405504366
$ws.Cells.Item(13, 17).Value = 0.05232025626733333
$ws.Cells.Item(13, 18).Value = 0.470882306406
$ws.Cells.Item(13, 19).Value = 0.0002645161645840909
$ws.Cells.Item(13, 20).Value = 0.0002645161645840909
